$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AD2").Value = 7.41818178187627
$ws.Range("AE2").Value = 7.399928440167592
$ws.Range("AD3").Value = -204943.4633497131
$ws.Range("AE3").Value = -201603.4482634706
